$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Hora da Movimentação" column right after "Data da Movimentação" ---
# Before this insert: Q=16 Valor da Unidade do Produto, R=17 Estoque Atual do Produto,
# S=19 Data da Movimentação, T=20 Motivo da Movimentação ...
# Inserting at column 20 (T) pushes "Motivo da Movimentação" and everything after it one
# column to the right, and the freshly inserted column becomes T.
$ws.Columns.Item(20).Insert()
$ws.Range("T1").Value = "Hora da Movimentação"

# --- Insert "Valor de Compra Unidade do Produto" column right before "Estoque Atual do Produto" ---
# Column 18 (R) is still "Estoque Atual do Produto" (untouched by the previous insert, since that
# happened further to the right). Inserting here pushes it (and everything after) one column right.
$ws.Columns.Item(18).Insert()
$ws.Range("R1").Value = "Valor de Compra Unidade do Produto"

# --- Column width adjustments that accompanied the new columns ---
$ws.Columns.Item(17).ColumnWidth = 29.83
$ws.Columns.Item(18).ColumnWidth = 34
$ws.Columns.Item(21).ColumnWidth = 23.17

# --- Update the active selection left on the sheet ---
$ws.Range("R2").Select()

Write-Output "edit applied"
